$d = $word.ActiveDocument

# --- Locate the paragraph "Diverts user away from bad traffic conditions when present"
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "^Diverts user away from bad traffic conditions when present") {
        $targetIndex = $i
        break
    }
}

$r = $d.Paragraphs.Item($targetIndex).Range
# The paragraph's range text is "...when present" followed by a manual line-break
# char (vertical tab) and the paragraph mark (carriage return), so the position
# right after the sentence - i.e. before the line break - is r.End - 2.
$insPos = $r.End - 2
$splitPoint = $d.Range($insPos, $insPos)
$splitPoint.InsertParagraphAfter()

# The split produced a brand-new paragraph (inheriting the same list formatting)
# that now holds only the manual line break. Insert the new bullet's text at its
# start so the break stays at the end of the paragraph, as in the original.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newStart = $newPara.Range.Start

$newText = "Automation of alerting user to avoidable traffic conditoins"
$textRange = $d.Range($newStart, $newStart)
$textRange.Text = $newText

# Re-anchor the "_GoBack" bookmark to sit right after the new text (and before the
# line break) in this new paragraph. Word only allows one bookmark per name, so
# this automatically removes it from its previous location (end of "Start Alerts").
$bmPos = $newStart + $newText.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
